# Update Fn1-Itga9 NATMI cell-cell communication TPM-derived metrics
# (ligand/receptor expression values and their derived specificity/edge weights)
# to match re-run of the scripts against new TPM data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 6.240107999999999
$ws.Cells.Item(2, 8).Value = 18.720324
$ws.Cells.Item(2, 9).Value = 0.01732230523539376
$ws.Cells.Item(2, 10).Value = 0.01732230523539376
$ws.Cells.Item(2, 13).Value = 0.7521946666666667
$ws.Cells.Item(2, 14).Value = 2.256584
$ws.Cells.Item(2, 15).Value = 0.07361670343069449
$ws.Cells.Item(2, 16).Value = 0.0736167034306945
$ws.Cells.Item(2, 17).Value = 4.693775957023999
$ws.Cells.Item(2, 18).Value = 42.243983613216
$ws.Cells.Item(2, 19).Value = 0.001275211007249949
$ws.Cells.Item(2, 20).Value = 0.001275211007249949
$ws.Cells.Item(3, 7).Value = 6.240107999999999
$ws.Cells.Item(3, 8).Value = 18.720324
$ws.Cells.Item(3, 9).Value = 0.01732230523539376
$ws.Cells.Item(3, 10).Value = 0.01732230523539376
$ws.Cells.Item(3, 15).Value = 0.6908862423022597
$ws.Cells.Item(3, 16).Value = 0.6908862423022598
$ws.Cells.Item(3, 17).Value = 44.050671682276
$ws.Cells.Item(3, 18).Value = 396.456045140484
$ws.Cells.Item(3, 19).Value = 0.01196774237209395
$ws.Cells.Item(3, 20).Value = 0.01196774237209396
$ws.Cells.Item(4, 7).Value = 6.240107999999999
$ws.Cells.Item(4, 8).Value = 18.720324
$ws.Cells.Item(4, 9).Value = 0.01732230523539376
$ws.Cells.Item(4, 10).Value = 0.01732230523539376
$ws.Cells.Item(4, 13).Value = 2.406242333333334
$ws.Cells.Item(4, 14).Value = 7.218727
$ws.Cells.Item(4, 15).Value = 0.2354970542670457
$ws.Cells.Item(4, 16).Value = 0.2354970542670457
$ws.Cells.Item(4, 17).Value = 15.015212034172
$ws.Cells.Item(4, 18).Value = 135.136908307548
$ws.Cells.Item(4, 19).Value = 0.004079351856049854
$ws.Cells.Item(4, 20).Value = 0.004079351856049854
$ws.Cells.Item(5, 9).Value = 0.9592798330716089
$ws.Cells.Item(5, 10).Value = 0.9592798330716091
$ws.Cells.Item(5, 13).Value = 0.7521946666666667
$ws.Cells.Item(5, 14).Value = 2.256584
$ws.Cells.Item(5, 15).Value = 0.07361670343069449
$ws.Cells.Item(5, 16).Value = 0.0736167034306945
$ws.Cells.Item(5, 17).Value = 259.9333377020453
$ws.Cells.Item(5, 18).Value = 2339.400039318408
$ws.Cells.Item(5, 19).Value = 0.07061901897827874
$ws.Cells.Item(5, 20).Value = 0.07061901897827877
$ws.Cells.Item(6, 9).Value = 0.9592798330716089
$ws.Cells.Item(6, 10).Value = 0.9592798330716091
$ws.Cells.Item(6, 15).Value = 0.6908862423022597
$ws.Cells.Item(6, 16).Value = 0.6908862423022598
$ws.Cells.Item(6, 19).Value = 0.6627532391871829
$ws.Cells.Item(6, 20).Value = 0.6627532391871831
$ws.Cells.Item(7, 9).Value = 0.9592798330716089
$ws.Cells.Item(7, 10).Value = 0.9592798330716091
$ws.Cells.Item(7, 13).Value = 2.406242333333334
$ws.Cells.Item(7, 14).Value = 7.218727
$ws.Cells.Item(7, 15).Value = 0.2354970542670457
$ws.Cells.Item(7, 16).Value = 0.2354970542670457
$ws.Cells.Item(7, 17).Value = 831.5169313749777
$ws.Cells.Item(7, 19).Value = 0.2259075749061472
$ws.Cells.Item(7, 20).Value = 0.2259075749061472
$ws.Cells.Item(8, 7).Value = 8.428738666666666
$ws.Cells.Item(8, 9).Value = 0.02339786169299727
$ws.Cells.Item(8, 10).Value = 0.02339786169299728
$ws.Cells.Item(8, 13).Value = 0.7521946666666667
$ws.Cells.Item(8, 14).Value = 2.256584
$ws.Cells.Item(8, 15).Value = 0.07361670343069449
$ws.Cells.Item(8, 16).Value = 0.0736167034306945
$ws.Cells.Item(8, 17).Value = 6.340052271793778
$ws.Cells.Item(8, 18).Value = 57.060470446144
$ws.Cells.Item(8, 19).Value = 0.001722473445165787
$ws.Cells.Item(8, 20).Value = 0.001722473445165788
$ws.Cells.Item(9, 7).Value = 8.428738666666666
$ws.Cells.Item(9, 9).Value = 0.02339786169299727
$ws.Cells.Item(9, 10).Value = 0.02339786169299728
$ws.Cells.Item(9, 15).Value = 0.6908862423022597
$ws.Cells.Item(9, 16).Value = 0.6908862423022598
$ws.Cells.Item(9, 17).Value = 59.50082910440622
$ws.Cells.Item(9, 18).Value = 535.507461939656
$ws.Cells.Item(9, 19).Value = 0.01616526074298287
$ws.Cells.Item(9, 20).Value = 0.01616526074298288
$ws.Cells.Item(10, 7).Value = 8.428738666666666
$ws.Cells.Item(10, 9).Value = 0.02339786169299727
$ws.Cells.Item(10, 10).Value = 0.02339786169299728
$ws.Cells.Item(10, 13).Value = 2.406242333333334
$ws.Cells.Item(10, 14).Value = 7.218727
$ws.Cells.Item(10, 15).Value = 0.2354970542670457
$ws.Cells.Item(10, 16).Value = 0.2354970542670457
$ws.Cells.Item(10, 19).Value = 0.005510127504848607
$ws.Cells.Item(10, 20).Value = 0.005510127504848609
